$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status column for the c29166f5 row (row 3) changes across the Overview
# sheet and both language sheets, since it's driven by the same underlying
# handback-transform result.
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Error Detail column (K) gets a new per-locale message describing the
# handback/handoff file name mismatch.
$zhcn.Range("K3").Value = "Handback file name: z5uzyeda.3kv is different with handoff file name: c29166f5-f96f-476a-8b00-752eaf6e6d6a.f22e27abdbcbda801205093e8840078375c197d3.zh-cn."
$dede.Range("K3").Value = "Handback file name: z5uzyeda.3kv is different with handoff file name: c29166f5-f96f-476a-8b00-752eaf6e6d6a.f22e27abdbcbda801205093e8840078375c197d3.de-de."
